$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Intro segment:" row/cell (A3) -- becomes the 2nd shared string
$ws.Range("A3").Value = "                              Intro segment:                                                              "

# Move the active selection to A4, matching the post-edit cursor position
$ws.Range("A4").Select()
